$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "datos actualizados" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 18:38"

# --- Re-rank countries: updated case counts move Cuba above Bulgaria
#     and Somalia above Niger/Andorra/Libano/Costa Rica (table sorted
#     descending by total cases, column B) ---
$ws.Range("A79").Value = "Cuba"
$ws.Range("A80").Value = "Bulgaria"
$ws.Range("A98").Value = "Somalia"
$ws.Range("A99").Value = "Niger"
$ws.Range("A100").Value = "Principado de Andorra"
$ws.Range("A101").Value = "Libano"
$ws.Range("A102").Value = "Costa Rica"

# --- Updated case/recovery/death statistics ---
# Row 4
$ws.Range("B4").Value = 1194434
$ws.Range("C4").Value = 6312
$ws.Range("D4").Value = 178934
$ws.Range("E4").Value = 946492
$ws.Range("G4").Value = 410
$ws.Range("H4").Value = 69008

# Row 6
$ws.Range("B6").Value = 211938
$ws.Range("C6").Value = 1221
$ws.Range("D6").Value = 82879
$ws.Range("E6").Value = 99980
$ws.Range("F6").Value = 1479
$ws.Range("G6").Value = 195
$ws.Range("H6").Value = 29079

# Row 7
$ws.Range("B7").Value = 190584
$ws.Range("C7").Value = 3985
$ws.Range("E7").Value = 161506
$ws.Range("G7").Value = 288
$ws.Range("H7").Value = 28734

# Row 15
$ws.Range("B15").Value = 59858
$ws.Range("C15").Value = 384
$ws.Range("D15").Value = 25422
$ws.Range("E15").Value = 30669
$ws.Range("G15").Value = 85
$ws.Range("H15").Value = 3767

# Row 47
$ws.Range("B47").Value = 7799
$ws.Range("C47").Value = 18
$ws.Range("D47").Value = 3786
$ws.Range("E47").Value = 3762
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 251

# Row 56
$ws.Range("B56").Value = 5053
$ws.Range("C56").Value = 150
$ws.Range("D56").Value = 1653
$ws.Range("E56").Value = 3221
$ws.Range("G56").Value = 5
$ws.Range("H56").Value = 179

# Row 70
$ws.Range("B70").Value = 2346
$ws.Range("C70").Value = 50
$ws.Range("D70").Value = 1544
$ws.Range("E70").Value = 704
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 98

# Row 71
$ws.Range("D71").Value = 1405
$ws.Range("E71").Value = 766

# Row 75
$ws.Range("B75").Value = 1984
$ws.Range("C75").Value = 52
$ws.Range("D75").Value = 1480
$ws.Range("E75").Value = 478
$ws.Range("F75").Value = 18
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 26

# Row 79
$ws.Range("B79").Value = 1668
$ws.Range("C79").Value = 19
$ws.Range("D79").Value = 876
$ws.Range("E79").Value = 723
$ws.Range("F79").Value = 9
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 69

# Row 80
$ws.Range("B80").Value = 1652
$ws.Range("C80").Value = 34
$ws.Range("D80").Value = 321
$ws.Range("E80").Value = 1253
$ws.Range("F80").Value = 40
$ws.Range("G80").Value = 5
$ws.Range("H80").Value = 78

# Row 98
$ws.Range("B98").Value = 756
$ws.Range("C98").Value = 34
$ws.Range("D98").Value = 61
$ws.Range("E98").Value = 660
$ws.Range("F98").Value = 2
$ws.Range("G98").Value = 3
$ws.Range("H98").Value = 35

# Row 99
$ws.Range("B99").Value = 750
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 518
$ws.Range("E99").Value = 196
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 36

# Row 100
$ws.Range("B100").Value = 748
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 493
$ws.Range("E100").Value = 210
$ws.Range("F100").Value = 17
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 45

# Row 101
$ws.Range("B101").Value = 740
$ws.Range("C101").Value = 3
$ws.Range("D101").Value = 200
$ws.Range("E101").Value = 515
$ws.Range("F101").Value = 43
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 25

# Row 102
$ws.Range("B102").Value = 739
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 386
$ws.Range("E102").Value = 347
$ws.Range("F102").Value = 6
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 6

